# Apply the commit's changes to the "Foglio1" worksheet of the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Ore Rendicontabili" (B) column values for rows 2-7, leaving the
# driving formulas in column C intact (they recalc to 0 since B is now blank).
$ws.Range("B2:B7").ClearContents()

# Update a few "Ore Totali" (D) figures which feed the Costo Committente (E) formulas.
$ws.Range("D4").Value = 250
$ws.Range("D6").Value = 150
$ws.Range("D7").Value = 220

# Move the active cell/selection to A2 (matches saved view state).
$ws.Range("A2").Select()

$wb.Save()
